# Weekly refresh of Femacal de La Calera - Bruselas (repollito) price rows.
# Re-assigns the Fecha/Volumen/Precio columns (D, J, K, L, M, P) across rows
# 2-32 to reflect the latest weekly data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44377
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 18000
$ws.Range("L2").Value = 19000
$ws.Range("M2").Value = 18500
$ws.Range("P2").Value = 1233

# Row 3
$ws.Range("D3").Value = 44326
$ws.Range("J3").Value = 45
$ws.Range("K3").Value = 15000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = 15000
$ws.Range("P3").Value = 1000

# Row 4
$ws.Range("D4").Value = 44327
$ws.Range("J4").Value = 35
$ws.Range("K4").Value = 15000
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = 15000
$ws.Range("P4").Value = 1000

# Row 5
$ws.Range("D5").Value = 44336
$ws.Range("J5").Value = 65
$ws.Range("K5").Value = 14000
$ws.Range("M5").Value = 14462
$ws.Range("P5").Value = 964

# Row 6
$ws.Range("D6").Value = 44333
$ws.Range("J6").Value = 35
$ws.Range("K6").Value = 15000
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = 15000
$ws.Range("P6").Value = 1000

# Row 7
$ws.Range("D7").Value = 44323
$ws.Range("J7").Value = 40

# Row 8
$ws.Range("D8").Value = 44308
$ws.Range("J8").Value = 40
$ws.Range("K8").Value = 16000
$ws.Range("L8").Value = 16000
$ws.Range("M8").Value = 16000
$ws.Range("P8").Value = 1067

# Row 10
$ws.Range("D10").Value = 44330
$ws.Range("J10").Value = 30
$ws.Range("K10").Value = 15000
$ws.Range("L10").Value = 15000
$ws.Range("M10").Value = 15000
$ws.Range("P10").Value = 1000

# Row 11
$ws.Range("D11").Value = 44328
$ws.Range("J11").Value = 38
$ws.Range("K11").Value = 15000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 15000
$ws.Range("P11").Value = 1000

# Row 12
$ws.Range("D12").Value = 44343

# Row 13
$ws.Range("D13").Value = 44344
$ws.Range("J13").Value = 40
$ws.Range("K13").Value = 20000
$ws.Range("L13").Value = 20000
$ws.Range("M13").Value = 20000
$ws.Range("P13").Value = 1333

# Row 14
$ws.Range("D14").Value = 44320
$ws.Range("J14").Value = 40

# Row 15
$ws.Range("D15").Value = 44406
$ws.Range("J15").Value = 50
$ws.Range("K15").Value = 22000
$ws.Range("L15").Value = 22000
$ws.Range("M15").Value = 22000
$ws.Range("P15").Value = 1467

# Row 16
$ws.Range("D16").Value = 44334
$ws.Range("J16").Value = 50
$ws.Range("K16").Value = 14000
$ws.Range("L16").Value = 14000
$ws.Range("M16").Value = 14000
$ws.Range("P16").Value = 933

# Row 17
$ws.Range("D17").Value = 44340
$ws.Range("J17").Value = 47
$ws.Range("K17").Value = 14000
$ws.Range("L17").Value = 14000
$ws.Range("M17").Value = 14000
$ws.Range("P17").Value = 933

# Row 18
$ws.Range("D18").Value = 44370
$ws.Range("J18").Value = 50
$ws.Range("K18").Value = 18000
$ws.Range("L18").Value = 18000
$ws.Range("M18").Value = 18000
$ws.Range("P18").Value = 1200

# Row 19
$ws.Range("D19").Value = 44309
$ws.Range("J19").Value = 50

# Row 20
$ws.Range("D20").Value = 44322
$ws.Range("J20").Value = 70
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = 14500
$ws.Range("P20").Value = 967

# Row 21
$ws.Range("D21").Value = 44397
$ws.Range("J21").Value = 73
$ws.Range("K21").Value = 21000
$ws.Range("L21").Value = 22000
$ws.Range("M21").Value = 21521
$ws.Range("P21").Value = 1435

# Row 22
$ws.Range("D22").Value = 44321
$ws.Range("J22").Value = 38
$ws.Range("K22").Value = 15000
$ws.Range("M22").Value = 15000
$ws.Range("P22").Value = 1000

# Row 23
$ws.Range("D23").Value = 44316
$ws.Range("J23").Value = 45
$ws.Range("K23").Value = 14000
$ws.Range("L23").Value = 15000
$ws.Range("M23").Value = 14444
$ws.Range("P23").Value = 963

# Row 24
$ws.Range("D24").Value = 44341
$ws.Range("J24").Value = 40

# Row 25
$ws.Range("D25").Value = 44315
$ws.Range("J25").Value = 65
$ws.Range("K25").Value = 14000
$ws.Range("M25").Value = 14538
$ws.Range("P25").Value = 969

# Row 26
$ws.Range("D26").Value = 44411
$ws.Range("J26").Value = 50
$ws.Range("K26").Value = 22000
$ws.Range("L26").Value = 22000
$ws.Range("M26").Value = 22000
$ws.Range("P26").Value = 1467

# Row 27
$ws.Range("D27").Value = 44313
$ws.Range("J27").Value = 40
$ws.Range("K27").Value = 14000
$ws.Range("L27").Value = 14000
$ws.Range("M27").Value = 14000
$ws.Range("P27").Value = 933

# Row 28
$ws.Range("D28").Value = 44329
$ws.Range("J28").Value = 35
$ws.Range("K28").Value = 15000
$ws.Range("M28").Value = 15000
$ws.Range("P28").Value = 1000

# Row 29
$ws.Range("D29").Value = 44438
$ws.Range("J29").Value = 75
$ws.Range("K29").Value = 19000
$ws.Range("L29").Value = 20000
$ws.Range("M29").Value = 19467
$ws.Range("P29").Value = 1298

# Row 30
$ws.Range("D30").Value = 44319
$ws.Range("J30").Value = 50

# Row 31
$ws.Range("D31").Value = 44314
$ws.Range("J31").Value = 45

# Row 32
$ws.Range("D32").Value = 44399
$ws.Range("J32").Value = 38
$ws.Range("K32").Value = 22000
$ws.Range("L32").Value = 22000
$ws.Range("M32").Value = 22000
$ws.Range("P32").Value = 1467
